$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / safe values (no numeric ambiguity) ---
$ws.Range("D2").Value = "63.424.49"
$ws.Range("E2").Value = "  -0.83%  "
$ws.Range("D3").Value = "3.239.74"
$ws.Range("E3").Value = "  +2.95%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  -1.42%  "
$ws.Range("E6").Value = "  -1.08%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.235.84"
$ws.Range("E8").Value = "  +3.03%  "
$ws.Range("E9").Value = "  -1.34%  "
$ws.Range("E10").Value = "  -1.09%  "
$ws.Range("E11").Value = "  -0.82%  "
$ws.Range("E13").Value = "  -2.52%  "
$ws.Range("E14").Value = "  -1.50%  "
$ws.Range("D15").Value = "3.770.57"
$ws.Range("E15").Value = "  +3.03%  "
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("D17").Value = "3.237.86"
$ws.Range("E17").Value = "  +3.00%  "
$ws.Range("D18").Value = "63.432.09"
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("E19").Value = "  -1.04%  "
$ws.Range("E20").Value = "  -2.46%  "
$ws.Range("E21").Value = "  -3.59%  "
$ws.Range("E22").Value = "  +2.30%  "
$ws.Range("E23").Value = "  +2.32%  "
$ws.Range("E24").Value = "  -4.72%  "
$ws.Range("E25").Value = "  -0.94%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  +8.04%  "
$ws.Range("E28").Value = "  -1.17%  "
$ws.Range("E29").Value = "  -1.33%  "
$ws.Range("E30").Value = "  +2.38%  "
$ws.Range("E31").Value = "  -0.44%  "
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("E33").Value = "  -3.34%  "
$ws.Range("E34").Value = "  -4.62%  "
$ws.Range("E35").Value = "  -1.66%  "
$ws.Range("E36").Value = "  -2.37%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").Value = "0.0₃0715"
$ws.Range("E38").Value = "  -4.38%  "
$ws.Range("E39").Value = "  -0.97%  "
$ws.Range("E41").Value = "  +0.73%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("E42").Value = "  -7.03%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.983.69"
$ws.Range("E43").Value = "  +1.58%  "
$ws.Range("E44").Value = "  -8.44%  "
$ws.Range("E45").Value = "  +2.50%  "
$ws.Range("E46").Value = "  -1.32%  "
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("E47").Value = "  -1.63%  "
$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("E49").Value = "  +0.34%  "
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("E51").Value = "  +11.36%  "

# --- Numeric-looking text values: force text via quote-prefix, then clear quote-prefix style ---
$ws.Range("D5").Value = "'594.70"
$ws.Range("D6").Value = "'142.13"
$ws.Range("D9").Value = "'0.520"
$ws.Range("D11").Value = "'5.34"
$ws.Range("D12").Value = "'0.466"
$ws.Range("D13").Value = "'0.0000248"
$ws.Range("D14").Value = "'34.46"
$ws.Range("D19").Value = "'6.80"
$ws.Range("D20").Value = "'475.88"
$ws.Range("D21").Value = "'14.16"
$ws.Range("D22").Value = "'0.728"
$ws.Range("D23").Value = "'7.90"
$ws.Range("D24").Value = "'84.02"
$ws.Range("D25").Value = "'13.18"
$ws.Range("D27").Value = "'7.57"
$ws.Range("D29").Value = "'8.09"
$ws.Range("D30").Value = "'2.11"
$ws.Range("D31").Value = "'27.50"
$ws.Range("D34").Value = "'2.53"
$ws.Range("D36").Value = "'5.92"
$ws.Range("D37").Value = "'52.75"
$ws.Range("D39").Value = "'0.0394"
$ws.Range("D40").Value = "'422.38"
$ws.Range("D41").Value = "'8.40"
$ws.Range("D42").Value = "'2.77"
$ws.Range("D45").Value = "'0.268"
$ws.Range("D46").Value = "'2.17"
$ws.Range("D47").Value = "'2.37"
$ws.Range("D48").Value = "'0.999"
$ws.Range("D49").Value = "'25.96"
$ws.Range("D51").Value = "'34.44"

# Reset style on the quote-prefixed cells so no stray formatting remains
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
